# Update the "想去人数" (interest count) figures for two sheets that mirror
# the same underlying event data: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6029
    $ws.Range("F4").Value = 180
    $ws.Range("F6").Value = 91
}
